# Re-order the comma-separated "Recorded By" names in column G into a
# canonical order: the literal "System" token always comes first, any
# other case-variant of "system" comes next, and everything else
# (actual user/email entries) follows in plain alphabetical order.
# Cells with only a single token are left untouched (no reordering is
# possible), and rows whose value is already canonical end up unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.Contains(",")) {
        $parts = $value -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        # Bucket 0: exact "System"; bucket 1: any other "system" casing;
        # bucket 2: everything else (sorted alphabetically within bucket).
        $exactSystem = @()
        $otherSystem = @()
        $others = @()
        foreach ($p in $trimmed) {
            if ($p.Equals("System")) {
                $exactSystem += $p
            } elseif ($p.ToLower() -eq "system") {
                $otherSystem += $p
            } else {
                $others += $p
            }
        }

        $others = $others | Sort-Object

        $ordered = @()
        $ordered += $exactSystem
        $ordered += $otherSystem
        $ordered += $others

        $newValue = [string]::Join(", ", $ordered)
        if ($newValue -ne $value) {
            $cell.Value2 = $newValue
        }
    }
}
